$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve textual formatting of the Price column (values like "1.002" or "28.576.12"
# must remain literal text, not be auto-converted to numbers/dates by Excel).
$priceRange = $ws.Range("D2:D51")
$origStyle = $ws.Range("D2").Style
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.570.24"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "1.922.41"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "315.95"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "0.5121"
$ws.Range("E7").Value = "  +1.97%  "
$ws.Range("D8").Value = "0.3987"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "0.09798"
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").Value = "1.149"
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").Value = "42.14"
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").Value = "6.484"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "21.04"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "1.912.75"
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("D15").Value = "7.426"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.00001135"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "94.32"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "0.06657"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").Value = "18.14"
$ws.Range("E20").Value = "  +4.03%  "
$ws.Range("D21").Value = "0.9996"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "6.291"
$ws.Range("E22").Value = "  +2.89%  "
$ws.Range("D23").Value = "28.628.44"
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("D24").Value = "11.49"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").Value = "2.310"
$ws.Range("E25").Value = "  +2.30%  "
$ws.Range("D26").Value = "2.703"
$ws.Range("E26").Value = "  +6.40%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.134.67"
$ws.Range("E27").Value = "  +2.24%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "21.25"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "158.25"
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "128.85"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "1.111"
$ws.Range("E31").Value = "  +4.28%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.1070"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.721"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.637"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "9.820"
$ws.Range("E35").Value = "  +4.32%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.06717"
$ws.Range("E36").Value = "  -1.75%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02447"
$ws.Range("E37").Value = "  +2.04%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "1.266"
$ws.Range("E38").Value = "  +3.79%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.2227"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "11.70"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.6456"
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").Value = "5.067"
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "1.189"
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.70"
$ws.Range("E45").Value = "  +1.87%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.6071"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.780"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "1.281"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.069"
$ws.Range("E49").Value = "  +3.57%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "124.46"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "1.200"
$ws.Range("E51").Value = "  -0.05%  "

# Restore original (default) style/number format now that text values are locked in.
$priceRange.Style = $origStyle

